# Apply the updates described by the commit diff to the "Report" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Numerical-variables correlation matrix (rows 46-48, cols B:D)
#    Only the off-diagonal cells changed value.
# ---------------------------------------------------------------------
$ws.Range("C46").Value = -0.005968836144085329
$ws.Range("D46").Value = 0.982716182458989

$ws.Range("B47").Value = -0.005968836144085329
$ws.Range("D47").Value = 0.01184095356250245

$ws.Range("B48").Value = 0.982716182458989
$ws.Range("C48").Value = 0.01184095356250245

# ---------------------------------------------------------------------
# 2) "Most frequent categorical values" table (rows 57-85).
#    Column A holds the "<Variable>_<value>" label, column B the share.
#    The whole block got regenerated (new random sample), so every label
#    from row 58 down is rewritten (row 57 keeps its original text) and
#    the frequencies are refreshed. Writing the labels in this exact
#    order reproduces the same shared-string allocation order as the
#    target workbook.
# ---------------------------------------------------------------------
$ws.Range("A58").Value = "Instruction Date_2014-02-04"
$ws.Range("A59").Value = "Instruction Date_2018-10-04"
$ws.Range("A60").Value = "Instruction Date_2014-02-26"
$ws.Range("A61").Value = "Country FI Initiation_CA"
$ws.Range("A62").Value = "Instruction Date_2014-07-14"
$ws.Range("A63").Value = "is_sender FI Initiation_Oui"
$ws.Range("A64").Value = "ind_lvts_Oui"
$ws.Range("A65").Value = "Country Sender_CA"
$ws.Range("A66").Value = "BIC FI Initiation_TDOMCATTTOR"
$ws.Range("A67").Value = "BIC FI Initiation_ROYCCAT2"
$ws.Range("A68").Value = "BIC FI Destination_ROYCCAT2"
$ws.Range("A69").Value = "BIC FI Initiation_CHASUS33FXR"
$ws.Range("A70").Value = "BIC Sender_TDOMCATTTOR"
$ws.Range("A71").Value = "BIC Sender_ROYCCAT2"
$ws.Range("A72").Value = "Channel_UNKNOWN"
$ws.Range("A73").Value = "Instruction Date_2016-02-15"
$ws.Range("A74").Value = "BIC FI Initiation_HKBCCATT"
$ws.Range("A75").Value = "BIC Sender_ROYCCAT3IMM"
$ws.Range("A76").Value = "BIC FI Initiation_ROYCCAT3IMM"
$ws.Range("A77").Value = "Instruction Date_2014-06-17"
$ws.Range("A78").Value = "BIC FI Initiation_BNDCCAMM"
$ws.Range("A79").Value = "BIC Sender_HKBCCATT"
$ws.Range("A80").Value = "Clients_Institutionnal"
$ws.Range("A81").Value = "Solution_Product_Direct_Debit"
$ws.Range("A82").Value = "Instruction_Withdrawal"
$ws.Range("A83").Value = "Instruction Date_2018-05-09"
$ws.Range("A84").Value = "Payment Engine_System 2"
$ws.Range("A85").Value = "Instruction Date_2014-07-31"

# Refreshed frequency values (column B) for the same block.
$ws.Range("B57").Value = 0.71
$ws.Range("B58").Value = 0.44
$ws.Range("B59").Value = 0.15
$ws.Range("B60").Value = 0.12
$ws.Range("B61").Value = 0.12
$ws.Range("B62").Value = 0.12
$ws.Range("B63").Value = 0.11
$ws.Range("B67").Value = 0.1
$ws.Range("B68").Value = 0.09
$ws.Range("B75").Value = 0.06
$ws.Range("B76").Value = 0.06
$ws.Range("B77").Value = 0.06
$ws.Range("B78").Value = 0.06
$ws.Range("B79").Value = 0.06
